$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.745293259620667
$ws.Range("B1").Value = 2.762935161590576
$ws.Range("C1").Value = 3.426137208938599
$ws.Range("D1").Value = 1.315690398216248
$ws.Range("E1").Value = 0.8754382133483887
